$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "o554F"

# Add new row 16 with data
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 1.010334526185207
$ws.Cells.Item(16, 4).Value = 0.926759645655817
$ws.Cells.Item(16, 5).Value = 1.009172679661215
$ws.Cells.Item(16, 6).Value = 1.010334526185207
$ws.Cells.Item(16, 7).Value = 0.9599991728739697
$ws.Cells.Item(16, 8).Value = 1.03691429481915
$ws.Cells.Item(16, 9).Value = 1.012069983939075
$ws.Cells.Item(16, 10).Value = 0.926759645655817
$ws.Cells.Item(16, 11).Value = 0.9679661626585159
$ws.Cells.Item(16, 12).Value = 0.9891503444218617
$ws.Cells.Item(16, 13).Value = 0.9925417171890726
